$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = "hi"
$ws.Range("F6").Select()
$ws.PageSetup.Orientation = 1
